# Rename header columns from *_old / *_new to *_FV2210 / *_FV2304,
# wrap the data range in an Excel Table (ListObject) and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update header row text -------------------------------------------------
$newHeadersFV2210 = @(
    "Segmentname_FV2210",
    "Segmentgruppe_FV2210",
    "Segment_FV2210",
    "Datenelement_FV2210",
    "Segment ID_FV2210",
    "Code_FV2210",
    "Qualifier_FV2210",
    "Beschreibung_FV2210",
    "Bedingungsausdruck_FV2210",
    "Bedingung_FV2210"
)

for ($i = 0; $i -lt $newHeadersFV2210.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $newHeadersFV2210[$i]
}

$newHeadersFV2304 = @(
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304"
)

# Column L is the 12th column (after "diff" in column K)
for ($i = 0; $i -lt $newHeadersFV2304.Length; $i++) {
    $ws.Cells.Item(1, 12 + $i).Value = $newHeadersFV2304[$i]
}

# --- 2. Turn the data range into an Excel Table (ListObject) -------------------
$tableRange = $ws.Range("A1:U55")
$listObject = $ws.ListObjects.Add(1, $tableRange, [System.Type]::Missing, 1)
$listObject.Name = "Table1"
$listObject.TableStyle = $null

# --- 3. Freeze the header row ---------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
